$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the individual card-detail rows (A2:A8) into a single Python-tuple-style
# string and place it in A2.
$name = $ws.Range("A2").Value()
$cost = $ws.Range("A3").Value()
$type = $ws.Range("A4").Value()
$ability1 = $ws.Range("A5").Value()
$ability2 = $ws.Range("A6").Value()
$ability3 = $ws.Range("A7").Value()
$pt = $ws.Range("A8").Value()

$combined = "('" + $name + "', ['" + $cost + "', '" + $type + "', '" + $ability1 + "', '" + $ability2 + "', '" + $ability3 + "', '" + $pt + "'])"

$ws.Range("A2").Value = $combined

# Remove the now-redundant rows (old A3:A8) so the used range collapses back
# to A1:A2.
$ws.Range("A3:A8").ClearContents() | Out-Null
